$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Date placeholder text: 12/12/2018 -> 11/26/2019 on the slide master, every
#    slide layout, and the notes master.
# ---------------------------------------------------------------------------
function Update-DatePlaceholder {
    param($shapes)
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            if ($sh.HasTextFrame) {
                $tr = $sh.TextFrame.TextRange
                if ($tr.Text -eq "12/12/2018") {
                    $tr.Text = "11/26/2019"
                }
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# The notes master's date placeholder doesn't take edits through its
# TextFrame directly in this host -- go through the HeadersFooters facade
# instead, which does persist for the notes master.
$notesMaster = $p.NotesMaster
$notesDateTime = $notesMaster.HeadersFooters.DateAndTime
$notesDateTime.Text = "11/26/2019"

# ---------------------------------------------------------------------------
# 2. Slide 1 (title slide): drop the old "TextBox 8" author list and add a
#    small "Updated Dec 2019" note under the existing credits text box.
# ---------------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
for ($i = $slide1.Shapes.Count; $i -ge 1; $i--) {
    $sh = $slide1.Shapes.Item($i)
    if ($sh.Name -eq "TextBox 8") {
        $sh.Delete() | Out-Null
    }
}

$note = $slide1.Shapes.AddTextbox(1, 46.79992125984252, 668.0109448818897, 163.8, 21.810944881889764)
$note.Name = "TextBox 5"
$noteTr = $note.TextFrame.TextRange
$noteTr.Text = "Updated Dec 2019"
$noteTr.Font.Size = 12
$noteTr.Font.Italic = $true
$note.TextFrame.WordWrap = $true
$note.TextFrame.AutoSize = 1
$note.Fill.Visible = 0
$note.Height = 21.810944881889764

# ---------------------------------------------------------------------------
# 3. Slide 3 (Linear regression and AIC): grow the "TextBox 6" shape a bit and
#    tweak the GLM bullet wording.
# ---------------------------------------------------------------------------
$slide3 = $p.Slides.Item(3)
for ($i = 1; $i -le $slide3.Shapes.Count; $i++) {
    $sh = $slide3.Shapes.Item($i)
    if ($sh.Name -eq "TextBox 6") {
        $sh.Height = 109.05472440944882
        $sh.TextFrame.TextRange.Replace("has some other distribution", "has some other statistical distribution") | Out-Null
    }
}
